$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the old "Value" column (E), shifting E:J -> F:K
$ws.Columns.Item(5).Insert()

# New column E header + left alignment for the whole Notes column
$ws.Range("E1").Value = "Notes"
$ws.Range("E1:E19").HorizontalAlignment = -4131
$ws.Columns.Item(5).ColumnWidth = 25.7109375

# Notes for the polarized-header rows
$ws.Range("E14").Value = "SWR25X"
$ws.Range("E15").Value = "SWH25X"
$ws.Range("E18").Value = "SWR25X"
$ws.Range("E19").Value = "SWH25X"

# New Digi-Key order part numbers in column D (Order PN) for header/socket rows
$ws.Range("D9").Value = "S7039-ND"
$ws.Range("D11").Value = "S7041-ND"
$ws.Range("D12").Value = "S9407-ND"
$ws.Range("D13").Value = "S9435-ND"
$ws.Range("D14").Value = "S9414-ND"
$ws.Range("D15").Value = "S9436-ND"
$ws.Range("D18").Value = "S9411-ND"
$ws.Range("D19").Value = "S9434-ND"

# Update the view/selection to match the authored edit
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D17").Select()

Write-Host "done"
